$d = $word.ActiveDocument

# 1. Change the body paragraph's style from "Author" to "Title"
$d.Paragraphs(1).Range.Style = "Title"

# 2. Title style: add spacing before/after (before=0/24pt -> 480/480 twips)
$titleStyle = $d.Styles("Title")
$titleStyle.ParagraphFormat.SpaceBefore = 24
$titleStyle.ParagraphFormat.SpaceAfter = 24

# 3. Author style: add new spacing before/after (18pt -> 360 twips each)
$authorStyle = $d.Styles("Author")
$authorStyle.ParagraphFormat.SpaceBefore = 18
$authorStyle.ParagraphFormat.SpaceAfter = 18

# 4. Date style: adjust spacing before/after (6pt before -> 120 twips, 16pt after -> 320 twips)
$dateStyle = $d.Styles("Date")
$dateStyle.ParagraphFormat.SpaceBefore = 6
$dateStyle.ParagraphFormat.SpaceAfter = 16
